$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3233.3333
$ws.Range("I40").Value = 4850
$ws.Range("J40").Value = 2425
$ws.Range("K40").Value = 4850
$ws.Range("L40").Value = 2425
$ws.Range("M40").Value = -4675
$ws.Range("N40").Value = -2775
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H70").Value = 176242
$ws.Range("I70").Value = 2500
$ws.Range("J70").Value = 349984
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 1049952
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -1050492
$ws.Range("H73").Value = 176242
$ws.Range("I73").Value = 2500
$ws.Range("J73").Value = 349984
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 1049952
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -1051824
$ws.Range("H105").Value = 43518.332
$ws.Range("J105").Value = 43518.332
$ws.Range("L105").Value = 43518.332
$ws.Range("N105").Value = -50506.332
$ws.Range("H106").Value = 71461220
$ws.Range("I106").Value = 76952456
$ws.Range("K106").Value = 76952456
$ws.Range("M106").Value = -76951825
$ws.Range("H131").Value = 1725
$ws.Range("I131").Value = 1725
$ws.Range("K131").Value = 5175
$ws.Range("M131").Value = -135

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20834460
$ws.Range("I2").Value = 23810110
$ws.Range("K2").Value = 23810110
$ws.Range("M2").Value = -23809997
$ws.Range("H61").Value = 2134.3333
$ws.Range("I61").Value = 2134.3333
$ws.Range("K61").Value = 2134.3333
$ws.Range("M61").Value = -1922.3333
$ws.Range("H116").Value = 20834460
$ws.Range("I116").Value = 23810110
$ws.Range("K116").Value = 23810110
$ws.Range("M116").Value = -23807816
$ws.Range("H122").Value = 928926.6
$ws.Range("I122").Value = 1269886.8
$ws.Range("K122").Value = 3809660.4
$ws.Range("M122").Value = -3807210.4
$ws.Range("H136").Value = 2134.3333
$ws.Range("I136").Value = 2134.3333
$ws.Range("K136").Value = 6402.999899999999
$ws.Range("M136").Value = -3852.999899999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20834460
$ws.Range("I3").Value = 23810110
$ws.Range("K3").Value = 23810110
$ws.Range("M3").Value = -23809996
$ws.Range("H80").Value = 391.18182
$ws.Range("I80").Value = 409.83334
$ws.Range("K80").Value = 409.83334
$ws.Range("M80").Value = 588.16666
$ws.Range("H83").Value = 391.18182
$ws.Range("I83").Value = 409.83334
$ws.Range("K83").Value = 2049.1667
$ws.Range("M83").Value = 2942.8333
$ws.Range("H92").Value = 24000
$ws.Range("J92").Value = 24000
$ws.Range("L92").Value = 24000
$ws.Range("N92").Value = -28992
$ws.Range("H105").Value = 4905576.5
$ws.Range("I105").Value = 8336930.5
$ws.Range("J105").Value = 3642.7144
$ws.Range("K105").Value = 8336930.5
$ws.Range("L105").Value = 3642.7144
$ws.Range("M105").Value = -8335183.5
$ws.Range("N105").Value = -7136.7144
$ws.Range("H107").Value = 1440.85
$ws.Range("I107").Value = 1407
$ws.Range("K107").Value = 1407
$ws.Range("M107").Value = 513

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 100000840
$ws.Range("I16").Value = 100000840
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 100000840
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -100000553
$ws.Range("H28").Value = 11888.667
$ws.Range("J28").Value = 12833
$ws.Range("L28").Value = 12833
$ws.Range("N28").Value = -13323
$ws.Range("H99").Value = 10055.968
$ws.Range("I99").Value = 6039.4707
$ws.Range("K99").Value = 6039.4707
$ws.Range("M99").Value = -4541.4707
$ws.Range("H105").Value = 2657
$ws.Range("I105").Value = 1191.909
$ws.Range("K105").Value = 1191.909
$ws.Range("M105").Value = 555.0909999999999
$ws.Range("H107").Value = 55556016
$ws.Range("I107").Value = 71428840
$ws.Range("K107").Value = 71428840
$ws.Range("M107").Value = -71426920
$ws.Range("H113").Value = 100000840
$ws.Range("I113").Value = 100000840
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 100000840
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -99998670
$ws.Range("H126").Value = 10055.968
$ws.Range("I126").Value = 6039.4707
$ws.Range("K126").Value = 18118.4121
$ws.Range("M126").Value = -15648.4121
$ws.Range("N16").ClearContents()
$ws.Range("N113").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83660.414
$ws.Range("J2").Value = 825.5
$ws.Range("L2").Value = 4953
$ws.Range("N2").Value = -5179
$ws.Range("H34").Value = 2478.5
$ws.Range("I34").Value = 609.6667
$ws.Range("J34").Value = 3599.8
$ws.Range("K34").Value = 1829.0001
$ws.Range("L34").Value = 10799.4
$ws.Range("M34").Value = -1745.0001
$ws.Range("N34").Value = -10967.4
$ws.Range("H38").Value = 71
$ws.Range("I38").Value = 79.5
$ws.Range("J38").Value = 59.666668
$ws.Range("K38").Value = 238.5
$ws.Range("L38").Value = 179.000004
$ws.Range("M38").Value = 108.5
$ws.Range("N38").Value = -873.000004
$ws.Range("H62").Value = 2999.8333
$ws.Range("J62").Value = 2999.8333
$ws.Range("L62").Value = 8999.499899999999
$ws.Range("N62").Value = -10371.4999
$ws.Range("H65").Value = 2999.8333
$ws.Range("J65").Value = 2999.8333
$ws.Range("L65").Value = 26998.4997
$ws.Range("N65").Value = -33862.4997
$ws.Range("H107").Value = 50625.15
$ws.Range("J107").Value = 67310.53
$ws.Range("L107").Value = 201931.59
$ws.Range("N107").Value = -205771.59

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 470
$ws.Range("I9").Value = 400
$ws.Range("J9").Value = 575
$ws.Range("K9").Value = 400
$ws.Range("L9").Value = 575
$ws.Range("M9").Value = -230
$ws.Range("N9").Value = -915
$ws.Range("H80").Value = 3488.5
$ws.Range("I80").Value = 2625.5715
$ws.Range("J80").Value = 5502
$ws.Range("K80").Value = 2625.5715
$ws.Range("L80").Value = 5502
$ws.Range("M80").Value = -1627.5715
$ws.Range("N80").Value = -7498
$ws.Range("H83").Value = 3488.5
$ws.Range("I83").Value = 2625.5715
$ws.Range("J83").Value = 5502
$ws.Range("K83").Value = 13127.8575
$ws.Range("L83").Value = 27510
$ws.Range("M83").Value = -8135.8575
$ws.Range("N83").Value = -37494
$ws.Range("H102").Value = 6686.9473
$ws.Range("I102").Value = 6780.6665
$ws.Range("K102").Value = 6780.6665
$ws.Range("M102").Value = -5158.6665
$ws.Range("H126").Value = 3995.6
$ws.Range("I126").Value = 3995.6
$ws.Range("K126").Value = 11986.8
$ws.Range("M126").Value = -9516.799999999999
$ws.Range("H132").Value = 1862.5454
$ws.Range("I132").Value = 1387.5555
$ws.Range("K132").Value = 4162.666499999999
$ws.Range("M132").Value = -1632.666499999999

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4457.091
$ws.Range("I22").Value = 4175.4287
$ws.Range("K22").Value = 4175.4287
$ws.Range("M22").Value = -3880.4287
$ws.Range("H27").Value = 4457.091
$ws.Range("I27").Value = 4175.4287
$ws.Range("K27").Value = 4175.4287
$ws.Range("M27").Value = -4068.4287
$ws.Range("H93").Value = 1284.2667
$ws.Range("I93").Value = 1076.4445
$ws.Range("J93").Value = 1596
$ws.Range("K93").Value = 1076.4445
$ws.Range("L93").Value = 1596
$ws.Range("M93").Value = 171.5554999999999
$ws.Range("N93").Value = -4092
$ws.Range("H134").Value = 57500
$ws.Range("J134").Value = 63000
$ws.Range("L134").Value = 63000
$ws.Range("N134").Value = -73140
$ws.Range("H136").Value = 8441.5
$ws.Range("I136").Value = 7995
$ws.Range("K136").Value = 23985
$ws.Range("M136").Value = -21435

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 3000
$ws.Range("I30").Value = 3000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2893
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982
$ws.Range("H113").Value = 863.3333
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 895
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 2685
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -7025
$ws.Range("H122").Value = 2400.2727
$ws.Range("I122").Value = 2289.2222
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 6867.6666
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -4417.6666
$ws.Range("N122").Value = -13600
$ws.Range("H136").Value = 1512.6666
$ws.Range("J136").Value = 3187
$ws.Range("L136").Value = 9561
$ws.Range("N136").Value = -14661
$ws.Range("N30").ClearContents()
